# Apply the sales-report update:
#  - MASS/Jan value (B4) drops from 3,750,000 to 2,400,000
#  - The BU total row (B2) is a SUM formula, so it recalculates automatically
#  - Move the active selection to D14 (matches the saved view state)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B4").Value = 2400000

$ws.Range("D14").Select()
